$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B ("ID Competição") holds the competition ID for rows 2..55.
# It was recorded as 56 but should be 256 (dropped leading digit).
for ($r = 2; $r -le 55; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    if ($cell.Value2 -eq 56) {
        $cell.Value = 256
    }
}
